$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '66.985.83'
$ws.Range("E2").Value = '  +1.17%  '

# Row 3
$ws.Range("D3").Value = '3.109.64'
$ws.Range("E3").Value = '  +2.94%  '

# Row 4
$ws.Range("E4").Value = '  -0.11%  '

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '579.91'
$ws.Range("E5").Value = '  +0.47%  '

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '172.94'
$ws.Range("E6").Value = '  +2.52%  '

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.999'
$ws.Range("E7").Value = '  -0.12%  '

# Row 8
$ws.Range("D8").Value = '3.104.30'
$ws.Range("E8").Value = '  +2.73%  '

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.524'
$ws.Range("E9").Value = '  +0.42%  '

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '6.43'
$ws.Range("E10").Value = '  -3.67%  '

# Row 11
$ws.Range("E11").Value = '  +1.22%  '

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.482'
$ws.Range("E12").Value = '  +0.04%  '

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.0000250'
$ws.Range("E13").Value = '  +0.79%  '

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '37.32'
$ws.Range("E14").Value = '  +0.21%  '

# Row 15
$ws.Range("E15").Value = '  +0.10%  '

# Row 16
$ws.Range("D16").Value = '3.618.84'
$ws.Range("E16").Value = '  +2.59%  '

# Row 17
$ws.Range("D17").Value = '66.933.92'
$ws.Range("E17").Value = '  +1.12%  '

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '7.21'
$ws.Range("E18").Value = '  -0.47%  '

# Row 19
$ws.Range("D19").Value = '3.103.70'
$ws.Range("E19").Value = '  +2.67%  '

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '16.34'
$ws.Range("E20").Value = '  +0.58%  '

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '485.78'
$ws.Range("E21").Value = '  +4.16%  '

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.721'
$ws.Range("E22").Value = '  +2.05%  '

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '7.58'
$ws.Range("E23").Value = '  +1.20%  '

# Row 24
$ws.Range("B24").Value = 'Litecoin'
$ws.Range("C24").Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '84.51'
$ws.Range("E24").Value = '  +1.26%  '

# Row 25
$ws.Range("B25").Value = 'InternetComputer(DFINITY)'
$ws.Range("C25").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '13.34'
$ws.Range("E25").Value = '  +2.86%  '

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '2.37'
$ws.Range("E26").Value = '  +3.42%  '

# Row 27
$ws.Range("B27").Value = 'RenderToken'
$ws.Range("C27").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '10.03'
$ws.Range("E27").Value = '  -1.52%  '

# Row 28
$ws.Range("B28").Value = 'Dai'
$ws.Range("C28").Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '1.00'
$ws.Range("E28").Value = '  -0.23%  '

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '8.09'
$ws.Range("E29").Value = '  -4.29%  '

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '2.43'
$ws.Range("E30").Value = '  -1.25%  '

# Row 31
$ws.Range("E31").Value = '  +1.81%  '

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '29.00'
$ws.Range("E32").Value = '  +2.72%  '

# Row 33
$ws.Range("E33").Value = '  -0.59%  '

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.115'
$ws.Range("E34").Value = '  -3.46%  '

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.999'
$ws.Range("E35").Value = '  -0.17%  '

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '5.94'
$ws.Range("E36").Value = '  +1.29%  '

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.990'
$ws.Range("E37").Value = '  +0.12%  '

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '47.47'
$ws.Range("E38").Value = '  -1.43%  '

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '2.12'
$ws.Range("E39").Value = '  +3.20%  '

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '50.19'
$ws.Range("E40").Value = '  +1.14%  '

# Row 41
$ws.Range("E41").Value = '  +1.33%  '

# Row 42
$ws.Range("E42").Value = '  +0.32%  '

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '8.68'
$ws.Range("E43").Value = '  +0.24%  '

# Row 44
$ws.Range("E44").Value = '  -2.51%  '

# Row 45
$ws.Range("D45").Value = '2.841.78'
$ws.Range("E45").Value = '  +3.93%  '

# Row 46
$ws.Range("B46").Value = 'Bittensor'
$ws.Range("C46").Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '387.70'
$ws.Range("E46").Value = '  +0.97%  '

# Row 47
$ws.Range("B47").Value = 'VeChain'
$ws.Range("C47").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.0359'
$ws.Range("E47").Value = '  -0.35%  '

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '136.08'
$ws.Range("E48").Value = '  +1.56%  '

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '24.97'
$ws.Range("E50").Value = '  -0.48%  '

# Row 51
$ws.Range("E51").Value = '  -1.03%  '
